# Reconstructs the source workbook's intended state per the commit diff:
#  - sheet renamed (new timestamp in the generated file name)
#  - last data row (the "...-Loading" / "加载中..." pair) removed, shrinking
#    the sheet from A1:C3 down to A1:C2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet to match the new export timestamp.
$ws.Name = "前端keyValue1696839541659.xlsx"

# 2) Remove row 3 (ce5b42-Loading / 加载中...) and shift the rows below it
#    (none, in this case) up — this also contracts the sheet's used range
#    / dimension from A1:C3 to A1:C2 automatically.
$ws.Rows.Item(3).Delete()
